$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# 1. Remove the stray _GoBack bookmark left over near the Base64 paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Re-language the six Message endpoint paragraphs from de-DE to en-US,
#    merging the split "Message"/" " runs in the [PUT] paragraph along the way.

$p96 = @"
<w:p $wNs>
<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>[POST]</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Message/</w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>GetMessagesInRoom</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>/</w:t></w:r>
</w:p>
"@
$d.Paragraphs.Item(96).Range.InsertXML($p96)

$p97 = @"
<w:p $wNs>
<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>[POST]</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Message/</w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>GetMessagesInRoomSince</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>/</w:t></w:r>
</w:p>
"@
$d.Paragraphs.Item(97).Range.InsertXML($p97)

$p98 = @"
<w:p $wNs>
<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>[POST]</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Message/</w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>GetNewMessagesForUser</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>/</w:t></w:r>
</w:p>
"@
$d.Paragraphs.Item(98).Range.InsertXML($p98)

$p99 = @"
<w:p $wNs>
<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>[DELETE]</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Message/</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Message</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>/</w:t></w:r>
</w:p>
"@
$d.Paragraphs.Item(99).Range.InsertXML($p99)

$p100 = @"
<w:p $wNs>
<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>[POST]</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Message/</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>Message </w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>/</w:t></w:r>
</w:p>
"@
$d.Paragraphs.Item(100).Range.InsertXML($p100)

$p101 = @"
<w:p $wNs>
<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>[PUT]</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Message/</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>Message </w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>/</w:t></w:r>
</w:p>
"@
$d.Paragraphs.Item(101).Range.InsertXML($p101)

# 3. AddFile signature now returns byte[] (the file's UUID) instead of void,
#    merge the split "byte[" / "] Content, " runs, and document the return
#    value with a trailing comment. Re-add a fresh _GoBack at the edit point.
$p120 = @"
<w:p $wNs>
<w:pPr><w:pStyle w:val='Heading4'/><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:tab/></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>byte[]</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>AddFile</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>(</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>byte[] Content, </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>C</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>File</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>CFile</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>);</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> // Returns UUID of file added</w:t></w:r>
</w:p>
"@
$d.Paragraphs.Item(120).Range.InsertXML($p120)

$pFile = $d.Paragraphs.Item(120)
$endRng = $d.Range($pFile.Range.End - 1, $pFile.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $endRng)
